$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2 = @{ B = 0.1169995834814548;  C = 0.002658071450198252; D = 0.7210945179870265; E = 0.5333859586016987;  G = 1.374138131520378 }
    3 = @{ B = 3.272327238179451;   C = 1.626987699542094;    D = 0.1496068669990043; E = 13.86384647080068;    G = 18.91276827552123 }
    4 = @{ B = 0.2881169905109251;  C = 1.626987699542094;    D = 0.7210945179870265; E = 0.5333859586016987;  G = 3.169585166641744 }
    5 = @{ B = 0.04172184405617529; C = 0.04103571897497393;  D = 0.1496068669990043; E = 0.5333859586016987;  G = 0.7657503886318522 }
    6 = @{ B = 0.1169995834814548;  C = 0.3048912486333797;   D = 0.1496068669990043; E = 0.5333859586016987;  G = 1.104883657715537 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
